# Recompute UTM easting/northing (Y_UTM / X_UTM), ZoneNumber and ZoneLetter
# for every data row, forcing the UTM zone number to 1 (the author's first
# attempt at re-running the lat/lon -> UTM conversion for the geological map,
# per the commit message). Latitude lives in column I ("Y"), longitude in
# column J ("X"); results are written back into columns B (Y_UTM), C (X_UTM),
# D (ZoneNumber) and E (ZoneLetter).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$wf = $excel.WorksheetFunction

# --- WGS84 / UTM constants (matches the standard `utm.from_latlon` forward
#     transverse-Mercator series used by the Python `utm` package) ---------
$K0 = 0.9996
$E  = 0.00669438
$EP2 = $E / (1.0 - $E)
$M1 = 1 - $E / 4 - 3 * $E * $E / 64 - 5 * $E * $E * $E / 256
$M2 = 3 * $E / 8 + 3 * $E * $E / 32 + 45 * $E * $E * $E / 1024
$M3 = 15 * $E * $E / 256 + 45 * $E * $E * $E / 1024
$M4 = 35 * $E * $E * $E / 3072
$R  = 6378137.0

# Zone is forced to 1 for every row -> central meridian = -177 degrees.
$forcedZone = 1
$centralLon = ($forcedZone - 1) * 6 - 180 + 3
$centralLonRad = $wf.Radians($centralLon)

$zoneLetters = "CDEFGHJKLMNPQRSTUVWXX"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 9).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 180 }

for ($r = 2; $r -le $lastRow; $r++) {
    $lat = $ws.Cells.Item($r, 9).Value()   # column I = Y (latitude)
    $lon = $ws.Cells.Item($r, 10).Value()  # column J = X (longitude)

    $latRad = $wf.Radians($lat)
    $latSin = $wf.Sin($latRad)
    $latCos = $wf.Cos($latRad)
    $latTan = $latSin / $latCos
    $latTan2 = $latTan * $latTan
    $latTan4 = $latTan2 * $latTan2

    $lonRad = $wf.Radians($lon)

    $n = $R / $wf.Sqrt(1 - $E * $latSin * $latSin)
    $c = $EP2 * $latCos * $latCos

    # longitude offset from the central meridian (always within +/-180 deg
    # for this dataset, so no extra angle-wrapping is required)
    $a = $latCos * ($lonRad - $centralLonRad)
    $a2 = $a * $a
    $a3 = $a2 * $a
    $a4 = $a3 * $a
    $a5 = $a4 * $a
    $a6 = $a5 * $a

    $m = $R * ($M1 * $latRad - $M2 * $wf.Sin(2 * $latRad) + $M3 * $wf.Sin(4 * $latRad) - $M4 * $wf.Sin(6 * $latRad))

    $easting = $K0 * $n * ($a + $a3 / 6 * (1 - $latTan2 + $c) + $a5 / 120 * (5 - 18 * $latTan2 + $latTan4 + 72 * $c - 58 * $EP2)) + 500000
    $northing = $K0 * ($m + $n * $latTan * ($a2 / 2 + $a4 / 24 * (5 - $latTan2 + 9 * $c + 4 * $c * $c) + $a6 / 720 * (61 - 58 * $latTan2 + $latTan4 + 600 * $c - 330 * $EP2)))

    if ($lat -lt 0) {
        $northing = $northing + 10000000
    }

    $letterIdx = $wf.Int(($lat + 80) / 8)
    $zoneLetter = $zoneLetters.Substring($letterIdx, 1)

    $ws.Cells.Item($r, 2).Value = $easting      # B: Y_UTM
    $ws.Cells.Item($r, 3).Value = $northing     # C: X_UTM
    $ws.Cells.Item($r, 4).Value = $forcedZone   # D: ZoneNumber
    $ws.Cells.Item($r, 5).Value = $zoneLetter   # E: ZoneLetter
}
